$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column K, copying formatting from J1 (same header style)
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "num_of_components"

# Populate num_of_components values for rows 2-87
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("K9").Value = 2
$ws.Range("K10").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("K12").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("K14").Value = 1
$ws.Range("K15").Value = 10
$ws.Range("K16").Value = 1
$ws.Range("K17").Value = 2
$ws.Range("K18").Value = 2
$ws.Range("K19").Value = 2
$ws.Range("K20").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("K22").Value = 1
$ws.Range("K23").Value = 2
$ws.Range("K24").Value = 1
$ws.Range("K25").Value = 1
$ws.Range("K26").Value = 1
$ws.Range("K27").Value = 2
$ws.Range("K28").Value = 1
$ws.Range("K29").Value = 1
$ws.Range("K30").Value = 1
$ws.Range("K31").Value = 3
$ws.Range("K32").Value = 1
$ws.Range("K33").Value = 1
$ws.Range("K34").Value = 1
$ws.Range("K35").Value = 1
$ws.Range("K36").Value = 1
$ws.Range("K37").Value = 1
$ws.Range("K38").Value = 2
$ws.Range("K39").Value = 1
$ws.Range("K40").Value = 1
$ws.Range("K41").Value = 1
$ws.Range("K42").Value = 1
$ws.Range("K43").Value = 1
$ws.Range("K44").Value = 1
$ws.Range("K45").Value = 1
$ws.Range("K46").Value = 1
$ws.Range("K47").Value = 1
$ws.Range("K48").Value = 1
$ws.Range("K49").Value = 1
$ws.Range("K50").Value = 1
$ws.Range("K51").Value = 1
$ws.Range("K52").Value = 1
$ws.Range("K53").Value = 1
$ws.Range("K54").Value = 1
$ws.Range("K55").Value = 1
$ws.Range("K56").Value = 1
$ws.Range("K57").Value = 1
$ws.Range("K58").Value = 1
$ws.Range("K59").Value = 1
$ws.Range("K60").Value = 1
$ws.Range("K61").Value = 1
$ws.Range("K62").Value = 1
$ws.Range("K63").Value = 1
$ws.Range("K64").Value = 1
$ws.Range("K65").Value = 1
$ws.Range("K66").Value = 1
$ws.Range("K67").Value = 1
$ws.Range("K68").Value = 1
$ws.Range("K69").Value = 1
$ws.Range("K70").Value = 1
$ws.Range("K71").Value = 1
$ws.Range("K72").Value = 1
$ws.Range("K73").Value = 3
$ws.Range("K74").Value = 1
$ws.Range("K75").Value = 1
$ws.Range("K76").Value = 1
$ws.Range("K77").Value = 1
$ws.Range("K78").Value = 1
$ws.Range("K79").Value = 1
$ws.Range("K80").Value = 1
$ws.Range("K81").Value = 1
$ws.Range("K82").Value = 1
$ws.Range("K83").Value = 1
$ws.Range("K84").Value = 1
$ws.Range("K85").Value = 1
$ws.Range("K86").Value = 2
$ws.Range("K87").Value = 1
